$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing numeric cells (rows 2-9, columns B-G) ---

# Row 2 (Q0)
$ws.Range("B2").Value = -0.03905235774959507
$ws.Range("C2").Value = 0.4387070408530346
$ws.Range("D2").Value = 0.45142052592237
$ws.Range("E2").Value = 0.6718783564919844
$ws.Range("F2").Value = 0.6774166724122599
$ws.Range("G2").Value = 51

# Row 3 (Q1)
$ws.Range("B3").Value = 0.0876619451090646
$ws.Range("C3").Value = 0.4671253201736847
$ws.Range("D3").Value = 0.4722669860490882
$ws.Range("E3").Value = 0.6872168406326261
$ws.Range("F3").Value = 0.6885227972073756
$ws.Range("G3").Value = 50

# Row 4 (Q2)
$ws.Range("B4").Value = -0.02200267238528682
$ws.Range("C4").Value = 0.4595796052308628
$ws.Range("D4").Value = 0.413274940495816
$ws.Range("E4").Value = 0.642864636214978
$ws.Range("F4").Value = 0.6491460788714786
$ws.Range("G4").Value = 49

# Row 5 (Q3)
$ws.Range("B5").Value = 0.1002432192375009
$ws.Range("C5").Value = 0.5011266898800307
$ws.Range("D5").Value = 0.4866691161701831
$ws.Range("E5").Value = 0.6976167401734158
$ws.Range("F5").Value = 0.6976827882510105
$ws.Range("G5").Value = 48

# Row 6 (Q4)
$ws.Range("B6").Value = 0.02274536467644276
$ws.Range("C6").Value = 0.4979849786295887
$ws.Range("D6").Value = 0.4567703583966128
$ws.Range("E6").Value = 0.6758478811068456
$ws.Range("F6").Value = 0.6827675668981734
$ws.Range("G6").Value = 47

# Row 7 (Q5)
$ws.Range("B7").Value = 0.092078373345108
$ws.Range("C7").Value = 0.5065592491403523
$ws.Range("D7").Value = 0.4958636947812035
$ws.Range("E7").Value = 0.7041758976145118
$ws.Range("F7").Value = 0.7058442120435615
$ws.Range("G7").Value = 46

# Row 8 (Q6)
$ws.Range("B8").Value = 0.03822042347955541
$ws.Range("C8").Value = 0.5402057006093556
$ws.Range("D8").Value = 0.4987337210178167
$ws.Range("E8").Value = 0.7062108191027781
$ws.Range("F8").Value = 0.7131441491908896
$ws.Range("G8").Value = 45

# Row 9 (Q7)
$ws.Range("B9").Value = 0.1352597487784972
$ws.Range("C9").Value = 0.5335606825057502
$ws.Range("D9").Value = 0.501733990090413
$ws.Range("E9").Value = 0.70833183614067
$ws.Range("F9").Value = 0.7033360171492812
$ws.Range("G9").Value = 44

# Row 10 (Q8) - also adds a new F10 cell that didn't exist before
$ws.Range("B10").Value = 0.05955804992731925
$ws.Range("C10").Value = 0.54934389505087
$ws.Range("D10").Value = 0.4945675297751405
$ws.Range("E10").Value = 0.7032549536086756
$ws.Range("F10").Value = 0.7090213886885923
$ws.Range("G10").Value = 43

# --- Add new row 11 (Q9) ---
$a11 = $ws.Range("A11")
$a11.Value = "Q9"
$a11.Font.Bold = $true
$a11.Borders.LineStyle = 1
$a11.HorizontalAlignment = -4108
$a11.VerticalAlignment = -4160

$ws.Range("B11").Value = 0.1423836843877604
$ws.Range("C11").Value = 0.5466645995897268
$ws.Range("D11").Value = 0.4954308102877077
$ws.Range("E11").Value = 0.7038684609269744
$ws.Range("F11").Value = 0.6976724939559396
$ws.Range("G11").Value = 42
